$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only rows 10 and 13-21 have textual content changes; rows 1-9, 11-12 keep their
# existing values untouched.

# Row 10: 'Objetivos:' | '5840942 - Marco Aurélio Kondracki de Alcântara' | '5840942 - Marco Aurélio Kondracki de Alcântara'
$ws.Range('A10').Value = 'Objetivos:'
$ws.Range('B10').Value = '5840942 - Marco Aurélio Kondracki de Alcântara'
$ws.Range('C10').Value = '5840942 - Marco Aurélio Kondracki de Alcântara'

# Row 13: 'Programa resumido:' | 'Semestral' | 'Semestral'
$ws.Range('A13').Value = 'Programa resumido:'
$ws.Range('B13').Value = 'Semestral'
$ws.Range('C13').Value = 'Semestral'

# Row 14: 'Short syllabus:' | '' | ''
$ws.Range('A14').Value = 'Short syllabus:'
$ws.Range('B14').Value = ''
$ws.Range('C14').Value = ''

# Row 15: 'Programa:' | '01/01/2012' | '01/01/2012'
$ws.Range('A15').Value = 'Programa:'
$ws.Range('B8').Copy($ws.Range('B15'))
$ws.Range('C8').Copy($ws.Range('C15'))

# Row 16: 'Syllabus:' | '' | ''
$ws.Range('A16').Value = 'Syllabus:'
$ws.Range('B16').Value = ''
$ws.Range('C16').Value = ''

# Row 17: 'Avaliação:' | '' | ''
$ws.Range('A17').Value = 'Avaliação:'
$ws.Range('B17').Value = ''
$ws.Range('C17').Value = ''

# Row 18: 'Método:' | '5840942 - Marco Aurélio Kondracki de Alcântara' | '5840942 - Marco Aurélio Kondracki de Alcântara'
$ws.Range('A18').Value = 'Método:'
$ws.Range('B19').Copy($ws.Range('B18'))
$ws.Range('B18').Value = '5840942 - Marco Aurélio Kondracki de Alcântara'
$ws.Range('C18').Value = '5840942 - Marco Aurélio Kondracki de Alcântara'

# Row 19: 'Critério:' | 'Aula expositiva e exercícios dirigidos.' | 'Aula expositiva e exercícios dirigidos.'
$ws.Range('A19').Value = 'Critério:'
$ws.Range('B19').Value = 'Aula expositiva e exercícios dirigidos.'
$ws.Range('C19').Value = 'Aula expositiva e exercícios dirigidos.'

# Row 20: 'Norma de recuperação:' | 'Média ponderada de exercícios e provas.' | 'Média ponderada de exercícios e provas.'
$ws.Range('A20').Value = 'Norma de recuperação:'
$ws.Range('B20').Value = 'Média ponderada de exercícios e provas.'
$ws.Range('C20').Value = 'Média ponderada de exercícios e provas.'

# Row 21: 'Bibliografia:' | 'Prova única com nota igual ou superior a 5,0.' | 'Prova única com nota igual ou superior a 5,0.'
$ws.Range('A21').Value = 'Bibliografia:'
$ws.Range('B21').Value = 'Prova única com nota igual ou superior a 5,0.'
$ws.Range('C21').Value = 'Prova única com nota igual ou superior a 5,0.'

# --- Adjust row heights that changed ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120

# --- Remove the now-obsolete last row (old long Bibliografia text row) ---
$ws.Rows.Item(22).Delete()

